# Auto-generated: apply scheduled market-data refresh to Asura_Profits workbook.
# Updates currentAveragePrice* / Leve Price* / LeveProfit* columns (H:N) for the
# rows whose underlying market data changed, per sheet (one sheet per job class).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4184.615
$ws.Range("I76").Value = 4600
$ws.Range("J76").Value = 3250
$ws.Range("K76").Value = 4600
$ws.Range("L76").Value = 3250
$ws.Range("M76").Value = -4285
$ws.Range("N76").Value = -3880
$ws.Range("H79").Value = 4184.615
$ws.Range("I79").Value = 4600
$ws.Range("J79").Value = 3250
$ws.Range("K79").Value = 4600
$ws.Range("L79").Value = 3250
$ws.Range("M79").Value = -3508
$ws.Range("N79").Value = -5434
$ws.Range("H92").Value = 691.8946999999999
$ws.Range("I92").Value = 730.06665
$ws.Range("J92").Value = 548.75
$ws.Range("K92").Value = 730.06665
$ws.Range("L92").Value = 548.75
$ws.Range("M92").Value = 517.93335
$ws.Range("N92").Value = -3044.75
$ws.Range("H100").Value = 2439
$ws.Range("I100").Value = 1761.6666
$ws.Range("K100").Value = 1761.6666
$ws.Range("M100").Value = -1220.6666
$ws.Range("H101").Value = 1400.4546
$ws.Range("J101").Value = 2685
$ws.Range("L101").Value = 8055
$ws.Range("N101").Value = -11299
$ws.Range("H113").Value = 3186.4285
$ws.Range("I113").Value = 2402.5
$ws.Range("K113").Value = 2402.5
$ws.Range("M113").Value = 851.5
$ws.Range("H116").Value = 11113497
$ws.Range("I116").Value = 50001400
$ws.Range("J116").Value = 2668.1428
$ws.Range("K116").Value = 50001400
$ws.Range("L116").Value = 2668.1428
$ws.Range("M116").Value = -49997958
$ws.Range("N116").Value = -9552.1428
$ws.Range("H125").Value = 4628.4165
$ws.Range("I125").Value = 6399.5
$ws.Range("J125").Value = 4274.2
$ws.Range("K125").Value = 57595.5
$ws.Range("L125").Value = 38467.8
$ws.Range("M125").Value = -55135.5
$ws.Range("N125").Value = -43387.8

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13683.518
$ws.Range("I32").Value = 15407.191
$ws.Range("J32").Value = 6318.727
$ws.Range("K32").Value = 15407.191
$ws.Range("L32").Value = 6318.727
$ws.Range("M32").Value = -15120.191
$ws.Range("N32").Value = -6892.727
$ws.Range("H74").Value = 1253
$ws.Range("I74").Value = 1360.6428
$ws.Range("K74").Value = 1360.6428
$ws.Range("M74").Value = -486.6428000000001
$ws.Range("H77").Value = 1253
$ws.Range("I77").Value = 1360.6428
$ws.Range("K77").Value = 6803.214
$ws.Range("M77").Value = -2435.214
$ws.Range("H88").Value = 2835
$ws.Range("J88").Value = 2943.75
$ws.Range("L88").Value = 2943.75
$ws.Range("N88").Value = -3755.75
$ws.Range("H91").Value = 2835
$ws.Range("J91").Value = 2943.75
$ws.Range("L91").Value = 2943.75
$ws.Range("N91").Value = -5751.75
$ws.Range("H97").Value = 665.4138
$ws.Range("I97").Value = 714.6087
$ws.Range("K97").Value = 714.6087
$ws.Range("M97").Value = -218.6087
$ws.Range("H104").Value = 25000
$ws.Range("J104").Value = 25000
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 52038
$ws.Range("I86").Value = 1729.8
$ws.Range("J86").Value = 102346.2
$ws.Range("K86").Value = 1729.8
$ws.Range("L86").Value = 102346.2
$ws.Range("M86").Value = -606.8
$ws.Range("N86").Value = -104592.2
$ws.Range("H89").Value = 52038
$ws.Range("I89").Value = 1729.8
$ws.Range("J89").Value = 102346.2
$ws.Range("K89").Value = 8649
$ws.Range("L89").Value = 511731
$ws.Range("M89").Value = -3033
$ws.Range("N89").Value = -522963

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1806.56
$ws.Range("I58").Value = 1433.7646
$ws.Range("J58").Value = 2598.75
$ws.Range("K58").Value = 1433.7646
$ws.Range("L58").Value = 2598.75
$ws.Range("M58").Value = -1230.7646
$ws.Range("N58").Value = -3004.75
$ws.Range("H105").Value = 2648.75
$ws.Range("I105").Value = 2648.75
$ws.Range("K105").Value = 2648.75
$ws.Range("M105").Value = -901.75
$ws.Range("H136").Value = 1806.56
$ws.Range("I136").Value = 1433.7646
$ws.Range("J136").Value = 2598.75
$ws.Range("K136").Value = 4301.293799999999
$ws.Range("L136").Value = 7796.25
$ws.Range("M136").Value = -1751.293799999999
$ws.Range("N136").Value = -12896.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 29921.5
$ws.Range("J53").Value = 29921.5
$ws.Range("L53").Value = 29921.5
$ws.Range("N53").Value = -31183.5
$ws.Range("H58").Value = 25772
$ws.Range("J58").Value = 24362.666
$ws.Range("L58").Value = 24362.666
$ws.Range("N58").Value = -24916.666
$ws.Range("H80").Value = 3525.5
$ws.Range("I80").Value = 3396.7856
$ws.Range("J80").Value = 3750.75
$ws.Range("K80").Value = 3396.7856
$ws.Range("L80").Value = 3750.75
$ws.Range("M80").Value = -2398.7856
$ws.Range("N80").Value = -5746.75
$ws.Range("H83").Value = 3525.5
$ws.Range("I83").Value = 3396.7856
$ws.Range("J83").Value = 3750.75
$ws.Range("K83").Value = 16983.928
$ws.Range("L83").Value = 18753.75
$ws.Range("M83").Value = -11991.928
$ws.Range("N83").Value = -28737.75
$ws.Range("H105").Value = 49950
$ws.Range("J105").Value = 49950
$ws.Range("L105").Value = 49950
$ws.Range("N105").Value = -56938
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1885.0938
$ws.Range("I93").Value = 1543.2106
$ws.Range("J93").Value = 2384.7693
$ws.Range("K93").Value = 1543.2106
$ws.Range("L93").Value = 2384.7693
$ws.Range("M93").Value = -295.2106000000001
$ws.Range("N93").Value = -4880.7693
$ws.Range("H104").Value = 21228.46
$ws.Range("J104").Value = 21228.46
$ws.Range("L104").Value = 21228.46
$ws.Range("N104").Value = -28216.46
$ws.Range("H106").Value = 40185
$ws.Range("J106").Value = 40185
$ws.Range("L106").Value = 40185
$ws.Range("N106").Value = -42709

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 69592.836
$ws.Range("I81").Value = 144473.62
$ws.Range("K81").Value = 288947.24
$ws.Range("M81").Value = -287886.24
$ws.Range("H84").Value = 69592.836
$ws.Range("I84").Value = 144473.62
$ws.Range("K84").Value = 1444736.2
$ws.Range("M84").Value = -1439432.2
$ws.Range("H104").Value = 267500
$ws.Range("J104").Value = 267500
$ws.Range("L104").Value = 267500
$ws.Range("N104").Value = -274488
$ws.Range("H113").Value = 1027.6
$ws.Range("I113").Value = 382.9091
$ws.Range("J113").Value = 2800.5
$ws.Range("K113").Value = 1148.7273
$ws.Range("L113").Value = 8401.5
$ws.Range("M113").Value = 1021.2727
$ws.Range("N113").Value = -12741.5
